$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 has no indel data at all (mSigHdp_ds_3k has no indel results) -
# make sure those cells are fully empty rather than empty-but-present.
$ws.Range("E2:G2").ClearContents()

# Fix the "Average indel" column (G): it was incorrectly set equal to
# indel_set1 (E) instead of being the average of indel_set1 (E) and
# indel_set2 (F).
$ws.Range("G3").Value = 4.7001681921296301
$ws.Range("G4").Value = 1.06960975694444
$ws.Range("G5").Value = 2.6826289212962999
$ws.Range("G6").Value = 1.0359959270833301

$ws.Range("G6").Select()
